# Edited shredding yield and Al content on the "Data" sheet (rows 14-17,
# column EW holds the end-of-series / 2050 value that the DT:EV columns
# linearly interpolate towards from column DS). Updating EW14:EW17 causes
# every shared formula in DT:EV on those rows to recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate() | Out-Null

$ws.Range("EW14").Value = 1.1000000000000001
$ws.Range("EW15").Value = 1
$ws.Range("EW16").Value = 0.8
$ws.Range("EW17").Value = 0.7

# Reflect the author's final on-screen selection for the sheet.
$ws.Range("EW31").Select() | Out-Null
